$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-11-08"

# Update the row label for the November partial-month row
$ws.Range("A13").Value = "November (through 11-08)"

# Row 13 - November (partial month) updated counts / rates
$ws.Range("C13").Value = 10
$ws.Range("F13").Value = 19
$ws.Range("I13").Value = 28
$ws.Range("J13").Value = 0.0345
$ws.Range("L13").Value = 16
$ws.Range("M13").Value = 0.2
$ws.Range("O13").Value = 11
$ws.Range("P13").Value = 0.0833
$ws.Range("R13").Value = 49
$ws.Range("S13").Value = 0.02
$ws.Range("U13").Value = 54
$ws.Range("V13").Value = 0.0182

# Row 14 - Totals updated to include new data
$ws.Range("C14").Value = 236
$ws.Range("D14").Value = 0.1194
$ws.Range("F14").Value = 453
$ws.Range("G14").Value = 0.103
$ws.Range("I14").Value = 677
$ws.Range("J14").Value = 0.0839
$ws.Range("L14").Value = 565
$ws.Range("M14").Value = 0.1102
$ws.Range("O14").Value = 445
$ws.Range("P14").Value = 0.0992
$ws.Range("R14").Value = 1052
$ws.Range("S14").Value = 0.0497
$ws.Range("U14").Value = 1414
$ws.Range("V14").Value = 0.0567
